$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the L&T rows (rows 28-41) from the bottom of the data range,
# which shrinks the sheet's used range from A1:K41 down to A1:K27.
$ws.Range("A28:K41").EntireRow.Delete()
